$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark column D data range as Text so purely numeric-looking
# price strings (e.g. "583.27") are kept as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '68.507.86'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '3.268.62'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '583.27'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '185.04'
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.599'
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').Value = '6.66'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = '0.420'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '3.850.11'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').Value = '28.31'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').Value = '68.524.43'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = '3.276.47'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '5.87'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '13.70'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = '394.56'
$ws.Range('E20').Value = '  +4.68%  '
$ws.Range('D21').Value = '7.71'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').Value = '71.63'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '0.520'
$ws.Range('E24').Value = '  +1.25%  '
$ws.Range('D25').Value = '0.0000120'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  +4.78%  '
$ws.Range('D27').Value = '9.78'
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').Value = '1.98'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').Value = '23.01'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').Value = '7.18'
$ws.Range('E32').Value = '  +3.11%  '
$ws.Range('D33').Value = '1.29'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('D36').Value = '163.43'
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('D37').Value = '1.97'
$ws.Range('E37').Value = '  +6.35%  '
$ws.Range('D38').Value = '0.825'
$ws.Range('E38').Value = '  -3.34%  '
$ws.Range('D39').Value = '26.73'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').Value = '4.60'
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('D41').Value = '6.57'
$ws.Range('E41').Value = '  -3.09%  '
$ws.Range('D42').Value = '2.53'
$ws.Range('E42').Value = '  -3.96%  '
$ws.Range('D43').Value = '25.60'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').Value = '0.0692'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('D45').Value = '41.35'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').Value = '2.657.42'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '339.36'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('E49').Value = '  +3.07%  '
$ws.Range('D50').Value = '31.78'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').Value = '0.996'
$ws.Range('E51').Value = '  -0.50%  '

# Restore normal styling on column D now that the text values are set.
$dRange.Style = "Normal"
